# Update cryptos list with latest prices and volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.141.72'
$ws.Range("E2").Value = '  -3.43%  '

$ws.Range("D3").Value = '3.299.03'
$ws.Range("E3").Value = '  -5.25%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '''543.97'
$ws.Range("E5").Value = '  -1.71%  '

$ws.Range("D6").Value = '''170.15'
$ws.Range("E6").Value = '  -4.60%  '

$ws.Range("E7").Value = '  -3.93%  '

$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").Value = '3.284.65'
$ws.Range("E9").Value = '  -5.53%  '

$ws.Range("D10").Value = '''0.605'
$ws.Range("E10").Value = '  -4.31%  '

$ws.Range("E11").Value = '  -1.21%  '

$ws.Range("D12").Value = '''52.61'
$ws.Range("E12").Value = '  -2.12%  '

$ws.Range("E13").Value = '  -2.64%  '

$ws.Range("D14").Value = '''8.78'
$ws.Range("E14").Value = '  -4.94%  '

$ws.Range("D15").Value = '3.832.74'
$ws.Range("E15").Value = '  -5.08%  '

$ws.Range("D16").Value = '''17.88'
$ws.Range("E16").Value = '  -4.54%  '

$ws.Range("D17").Value = '3.307.05'
$ws.Range("E17").Value = '  -5.03%  '

$ws.Range("E18").Value = '  -4.17%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '''11.57'
$ws.Range("E19").Value = '  -4.02%  '

$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D20").Value = '63.232.45'
$ws.Range("E20").Value = '  -3.32%  '

$ws.Range("D21").Value = '''0.965'
$ws.Range("E21").Value = '  -2.42%  '

$ws.Range("D22").Value = '''411.85'
$ws.Range("E22").Value = '  -1.24%  '

$ws.Range("E23").Value = '  -0.60%  '

$ws.Range("E24").Value = '  +6.02%  '

$ws.Range("D25").Value = '''13.58'
$ws.Range("E25").Value = '  +4.71%  '

$ws.Range("D26").Value = '''82.38'
$ws.Range("E26").Value = '  -4.43%  '

$ws.Range("E27").Value = '  -2.94%  '

$ws.Range("D28").Value = '''2.69'
$ws.Range("E28").Value = '  -5.19%  '

$ws.Range("E29").Value = '  -5.98%  '

$ws.Range("E30").Value = '  -4.86%  '

$ws.Range("D31").Value = '''6.31'
$ws.Range("E31").Value = '  -3.38%  '

$ws.Range("D32").Value = '''11.25'
$ws.Range("E32").Value = '  -4.14%  '

$ws.Range("D33").Value = '''564.73'
$ws.Range("E33").Value = '  -8.09%  '

$ws.Range("E34").Value = '  -3.91%  '

$ws.Range("D35").Value = '''57.43'
$ws.Range("E35").Value = '  -3.36%  '

$ws.Range("D36").Value = '''0.999'
$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("E37").Value = '  -0.60%  '

$ws.Range("E38").Value = '  -6.86%  '

$ws.Range("D39").Value = '''3.38'
$ws.Range("E39").Value = '  +4.51%  '

$ws.Range("D40").Value = '0.0₃0730'
$ws.Range("E40").Value = '  -6.80%  '

$ws.Range("D41").Value = '''0.361'
$ws.Range("E41").Value = '  -4.70%  '

$ws.Range("D42").Value = '3.106.12'
$ws.Range("E42").Value = '  -7.80%  '

$ws.Range("D43").Value = '''0.999'
$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("E44").Value = '  -2.52%  '

$ws.Range("D45").Value = '''3.20'
$ws.Range("E45").Value = '  -1.62%  '

$ws.Range("E47").Value = '  -5.43%  '

$ws.Range("E48").Value = '  -4.21%  '

$ws.Range("E49").Value = '  -3.92%  '

$ws.Range("D50").Value = '''132.25'
$ws.Range("E50").Value = '  -3.82%  '

$ws.Range("D51").Value = '''7.95'
$ws.Range("E51").Value = '  -6.27%  '
